# Applies the "updated 4.0 files and mdl" changes:
#  - About sheet: bump the "last updated" date in C1 from 45320 to 45392
#  - MCF sheet: raise several capacity-factor values from 0.85/0.95 to 1
#  - MCF sheet: move the active cell selection to B17

$wb = $excel.ActiveWorkbook

# --- About sheet: update date stamp in C1 ---
$wsAbout = $wb.Worksheets.Item("About")
$wsAbout.Range("C1").Value = 45392

# --- MCF sheet: update capacity factor values ---
$wsMcf = $wb.Worksheets.Item("MCF")

$rows = @(2, 3, 4, 6, 10, 11, 12, 13, 14, 16, 17, 18)
foreach ($r in $rows) {
    $wsMcf.Cells.Item($r, 2).Value = 1
}

# Move/activate the sheet and set the selection to match the new cursor position
$wsMcf.Activate()
$wsMcf.Range("B17").Select()
